$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "86×91=" "47×91="
Replace-Text "40×17=" "33×11="
Replace-Text "39×57=" "67×69="
Replace-Text "15×76=" "55×15="
Replace-Text "21×48=" "21×78="
Replace-Text "18×43=" "56×71="
Replace-Text "51×79=" "73×96="
Replace-Text "25×62=" "82×14="
Replace-Text "44×74=" "27×99="
Replace-Text "47×85=" "23×12="
Replace-Text "92×26=" "84×73="
Replace-Text "19×49=" "15×38="
Replace-Text "12×61=" "20×42="
Replace-Text "71×25=" "67×51="
Replace-Text "66×16=" "91×74="
Replace-Text "40×86=" "14×48="
Replace-Text "74×75=" "27×26="
Replace-Text "48×37=" "87×32="
Replace-Text "15×17=" "42×32="
Replace-Text "42×19=" "36×39="
Replace-Text "20×37=" "17×49="
Replace-Text "94×15=" "87×70="
Replace-Text "15×99=" "60×85="
Replace-Text "43×33=" "38×86="
Replace-Text "82×87=" "71×92="
